$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume/1h change (E) columns for rows with new crypto data.
# D-column values that look like plain numbers need the cell pre-formatted as Text
# ("@") before assignment, otherwise Excel auto-converts them to numeric values
# (losing the decimal-grouped "thousands" formatting / introducing float rounding).
$ws.Range("D2").Value = "47.562.35"
$ws.Range("E2").Value = "  +5.52%  "
$ws.Range("D3").Value = "2.497.84"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.58"
$ws.Range("E5").Value = "  +2.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.14"
$ws.Range("E6").Value = "  +2.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.522"
$ws.Range("E7").Value = "  +1.84%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("E9").Value = "  +2.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.75"
$ws.Range("E10").Value = "  +6.62%  "
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.14"
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("D15").Value = "2.884.22"
$ws.Range("E15").Value = "  +2.51%  "
$ws.Range("D16").Value = "2.520.44"
$ws.Range("E16").Value = "  +3.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.842"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "47.417.30"
$ws.Range("E18").Value = "  +5.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.67"
$ws.Range("E19").Value = "  +3.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.55"
$ws.Range("E20").Value = "  +2.97%  "
$ws.Range("D21").Value = "0.0₃0933"
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.76"
$ws.Range("E22").Value = "  +2.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.86"
$ws.Range("E23").Value = "  +2.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  +5.65%  "
$ws.Range("E25").Value = "  +3.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.21"
$ws.Range("E26").Value = "  +3.73%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.08"
$ws.Range("E28").Value = "  +5.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  -3.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.21"
$ws.Range("E30").Value = "  +7.54%  "
$ws.Range("E31").Value = "  +7.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.50"
$ws.Range("E32").Value = "  +0.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.95"
$ws.Range("E33").Value = "  -1.18%  "
$ws.Range("E34").Value = "  +3.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0782"
$ws.Range("E35").Value = "  +2.67%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  +4.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.65"
$ws.Range("E38").Value = "  +5.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.99"
$ws.Range("E39").Value = "  +3.99%  "
$ws.Range("E40").Value = "  +2.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.20"
$ws.Range("E43").Value = "  +1.89%  "
$ws.Range("E44").Value = "  +2.40%  "
$ws.Range("D45").Value = "1.970.77"
$ws.Range("E45").Value = "  +1.91%  "
$ws.Range("E46").Value = "  +1.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.10"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("E48").Value = "  +1.85%  "
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.34"
$ws.Range("E50").Value = "  +13.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.83"
$ws.Range("E51").Value = "  +3.14%  "

# Rows 41 and 42 swap coin identity (Monero <-> WEMIXToken) along with new values.
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.23"
$ws.Range("E41").Value = "  +0.42%  "

$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "121.04"
$ws.Range("E42").Value = "  -3.40%  "
